$d = $word.ActiveDocument

$pairs = @(
    @{old="58+17=75"; new="97-94=3"},
    @{old="42-11=31"; new="76-18=58"},
    @{old="77+21=98"; new="89-33=56"},
    @{old="70+21=91"; new="19+10=29"},
    @{old="43-37=6"; new="76-59=17"},
    @{old="86-29=57"; new="79-68=11"},
    @{old="45+37=82"; new="51-4=47"},
    @{old="35-1=34"; new="58-26=32"},
    @{old="86-12=74"; new="37-30=7"},
    @{old="55+39=94"; new="39+57=96"},
    @{old="30+6=36"; new="87-45=42"},
    @{old="3-2=1"; new="51-36=15"},
    @{old="44-8=36"; new="20+63=83"},
    @{old="30-26=4"; new="28+21=49"},
    @{old="11+25=36"; new="85-52=33"},
    @{old="15+24=39"; new="79-63=16"},
    @{old="69-19=50"; new="93-85=8"},
    @{old="45+30=75"; new="56-3=53"},
    @{old="75+18=93"; new="11+67=78"},
    @{old="58+20=78"; new="29+5=34"},
    @{old="75+15=90"; new="90-70=20"},
    @{old="22+25=47"; new="68-42=26"},
    @{old="35+64=99"; new="7+6=13"},
    @{old="77-13=64"; new="15+66=81"},
    @{old="75-60=15"; new="12-4=8"},
    @{old="48-38=10"; new="85-35=50"},
    @{old="96-2=94"; new="2+70=72"},
    @{old="66-5=61"; new="72-72=0"},
    @{old="31-23=8"; new="30-8=22"},
    @{old="34+39=73"; new="95-34=61"},
    @{old="81-67=14"; new="86-58=28"},
    @{old="53+35=88"; new="85-17=68"},
    @{old="30+19=49"; new="70-41=29"},
    @{old="57-54=3"; new="71+16=87"},
    @{old="18-17=1"; new="57-41=16"},
    @{old="14+64=78"; new="36+49=85"},
    @{old="25+45=70"; new="91-17=74"},
    @{old="36+14=50"; new="47-1=46"},
    @{old="72-40=32"; new="28+4=32"},
    @{old="88-42=46"; new="75-26=49"},
    @{old="36-24=12"; new="61+37=98"},
    @{old="64-45=19"; new="57+9=66"},
    @{old="9-4=5"; new="10+32=42"},
    @{old="91-78=13"; new="11+13=24"},
    @{old="58-12=46"; new="27+13=40"},
    @{old="64-57=7"; new="24+75=99"},
    @{old="47+12=59"; new="2+54=56"},
    @{old="87-58=29"; new="14+56=70"},
    @{old="16+29=45"; new="62-23=39"},
    @{old="66-38=28"; new="81-43=38"},
    @{old="17-7=10"; new="60+13=73"},
    @{old="55-26=29"; new="75-52=23"},
    @{old="42-31=11"; new="1+93=94"},
    @{old="58+19=77"; new="30-13=17"},
    @{old="89-48=41"; new="13+0=13"},
    @{old="87-2=85"; new="66-41=25"},
    @{old="52-24=28"; new="59+31=90"},
    @{old="93-73=20"; new="49-17=32"},
    @{old="69-23=46"; new="39+6=45"},
    @{old="45-38=7"; new="67+7=74"},
    @{old="82-43=39"; new="35-13=22"},
    @{old="42-18=24"; new="2+4=6"},
    @{old="23-16=7"; new="37+20=57"},
    @{old="88-41=47"; new="59-51=8"},
    @{old="20+26=46"; new="82-28=54"},
    @{old="42-38=4"; new="1+83=84"},
    @{old="29+70=99"; new="51-25=26"},
    @{old="28+53=81"; new="17+23=40"},
    @{old="78+21=99"; new="75-52=23"},
    @{old="89-30=59"; new="7+76=83"},
    @{old="92-0=92"; new="90-12=78"},
    @{old="54-21=33"; new="82-21=61"},
    @{old="72-40=32"; new="15+17=32"},
    @{old="95-83=12"; new="11+26=37"},
    @{old="38+35=73"; new="68-31=37"},
    @{old="29+50=79"; new="6+1=7"},
    @{old="28+54=82"; new="76-24=52"},
    @{old="91-45=46"; new="7-1=6"},
    @{old="87-66=21"; new="84+2=86"},
    @{old="16+56=72"; new="21+14=35"},
    @{old="6+27=33"; new="84-57=27"},
    @{old="61-39=22"; new="40-12=28"},
    @{old="8+61=69"; new="82-13=69"},
    @{old="19+23=42"; new="2+63=65"},
    @{old="55-37=18"; new="54+28=82"},
    @{old="40-16=24"; new="37+41=78"},
    @{old="33+13=46"; new="35+19=54"},
    @{old="59-2=57"; new="22+0=22"},
    @{old="97-83=14"; new="25+60=85"},
    @{old="13+80=93"; new="53+8=61"},
    @{old="28+66=94"; new="24+1=25"},
    @{old="7+63=70"; new="40-2=38"},
    @{old="57+37=94"; new="27+20=47"},
    @{old="95-69=26"; new="2+60=62"},
    @{old="17+60=77"; new="10+11=21"},
    @{old="61-29=32"; new="17+66=83"},
    @{old="20-19=1"; new="75+0=75"},
    @{old="90-31=59"; new="20+37=57"},
    @{old="58+2=60"; new="49-35=14"},
    @{old="9+38=47"; new="50-41=9"}
)

$cursorStart = 0
foreach ($p in $pairs) {
    $rng = $d.Range($cursorStart, $d.Content.End)
    $rng.Find.Execute($p.old, $false, $false, $false, $false, $false, $true, 1, $false, $p.new, 1)
    $cursorStart = $rng.End
}
